$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9048283336185818
$ws.Range("C2").Value = 0.1355680638631682
$ws.Range("D2").Value = 0.1525889375752456
$ws.Range("F2").Value = 1.902693688906858
$ws.Range("G2").Value = 0.002505473217499012
$ws.Range("J2").Value = 0.2234386864452205
$ws.Range("K2").Value = 0.4468790527791953
$ws.Range("L2").Value = 0.3303546847548375
$ws.Range("N2").Value = 2.201099405360304
$ws.Range("O2").Value = 4.893947308165821
$ws.Range("B3").Value = 0.8668321265355701
$ws.Range("C3").Value = 0.1350393194915753
$ws.Range("D3").Value = 0.1505431398048316
$ws.Range("F3").Value = 1.908835776594373
$ws.Range("G3").Value = 0.00250794563295973
$ws.Range("J3").Value = 0.2244571156854427
$ws.Range("K3").Value = 0.4130303882693624
$ws.Range("L3").Value = 0.3255741744138874
$ws.Range("N3").Value = 2.221283881089278
$ws.Range("O3").Value = 4.917287170159483
$ws.Range("B4").Value = 0.8438494580375391
$ws.Range("C4").Value = 0.1347176615995629
$ws.Range("D4").Value = 0.1493430115258221
$ws.Range("F4").Value = 1.913412334772815
$ws.Range("G4").Value = 0.00250954619230731
$ws.Range("J4").Value = 0.2251567073599006
$ws.Range("K4").Value = 0.3923674299748541
$ws.Range("L4").Value = 0.3227714864182047
$ws.Range("N4").Value = 2.234307760467949
$ws.Range("O4").Value = 4.933709342383509
$ws.Range("B5").Value = 0.8345718804397109
$ws.Range("C5").Value = 0.1345873525378316
$ws.Range("D5").Value = 0.1488681060762715
$ws.Range("F5").Value = 1.915480059843929
$ws.Range("G5").Value = 0.002510219235759086
$ws.Range("J5").Value = 0.2254605035722577
$ws.Range("K5").Value = 0.3839778903645765
$ws.Range("L5").Value = 0.3216628456813027
$ws.Range("N5").Value = 2.239773701771396
$ws.Range("O5").Value = 4.94092785882394
$ws.Range("B6").Value = 0.8330366856419289
$ws.Range("C6").Value = 0.134565761756285
$ws.Range("D6").Value = 0.1487901056145731
$ws.Range("F6").Value = 1.915835656207335
$ws.Range("G6").Value = 0.002510332252328501
$ws.Range("J6").Value = 0.2255120794557897
$ws.Range("K6").Value = 0.3825866889657448
$ws.Range("L6").Value = 0.3214807828869297
$ws.Range("N6").Value = 2.24069089808975
$ws.Range("O6").Value = 4.942158291302974
$ws.Range("B7").Value = 0.8437239798888356
$ws.Range("C7").Value = 0.1347159010693701
$ws.Range("D7").Value = 0.1493365493650529
$ws.Range("F7").Value = 1.913439399625389
$ws.Range("G7").Value = 0.002509555184877985
$ws.Range("J7").Value = 0.2251607286774657
$ws.Range("K7").Value = 0.3922541603471359
$ws.Range("L7").Value = 0.3227563991628841
$ws.Range("N7").Value = 2.234380833809312
$ws.Range("O7").Value = 4.933804562066015
$ws.Range("B8").Value = 0.8916556263171174
$ws.Range("C8").Value = 0.135385141215238
$ws.Range("D8").Value = 0.1518719677800817
$ws.Range("F8").Value = 1.904644485968369
$ws.Range("G8").Value = 0.002506308622634304
$ws.Range("J8").Value = 0.2237744439504965
$ws.Range("K8").Value = 0.4351833888030114
$ws.Range("L8").Value = 0.3286789345736452
$ws.Range("N8").Value = 2.207928179819423
$ws.Range("O8").Value = 4.901561086116942
$ws.Range("B9").Value = 0.9883763353728625
$ws.Range("C9").Value = 0.1367206189942181
$ws.Range("D9").Value = 0.1572853265539891
$ws.Range("F9").Value = 1.893777787850354
$ws.Range("G9").Value = 0.002500593846905722
$ws.Range("J9").Value = 0.2216439648287185
$ws.Range("K9").Value = 0.5203021391653522
$ws.Range("L9").Value = 0.3413394481077887
$ws.Range("N9").Value = 2.161054420375386
$ws.Range("O9").Value = 4.854908595343574
$ws.Range("B10").Value = 1.061070200840447
$ws.Range("C10").Value = 0.1377151089274093
$ws.Range("D10").Value = 0.1615282348475517
$ws.Range("F10").Value = 1.889672423477563
$ws.Range("G10").Value = 0.002496788630525896
$ws.Range("J10").Value = 0.2204355115556176
$ws.Range("K10").Value = 0.5833888419926154
$ws.Range("L10").Value = 0.351273097394639
$ws.Range("N10").Value = 2.129658820328546
$ws.Range("O10").Value = 4.830718478016479
$ws.Range("B11").Value = 1.094489421091566
$ws.Range("C11").Value = 0.1381702638518263
$ws.Range("D11").Value = 0.163515456653343
$ws.Range("F11").Value = 1.888644727037502
$ws.Range("G11").Value = 0.002495142144328953
$ws.Range("J11").Value = 0.2199628878388467
$ws.Range("K11").Value = 0.6122042051141818
$ws.Range("L11").Value = 0.3559282106799913
$ws.Range("N11").Value = 2.116035643752155
$ws.Range("O11").Value = 4.821899726791088
$ws.Range("B12").Value = 1.107194141124012
$ws.Range("C12").Value = 0.1383430003959703
$ws.Range("D12").Value = 0.1642761080612161
$ws.Range("F12").Value = 1.888376125124864
$ws.Range("G12").Value = 0.002494530755271285
$ws.Range("J12").Value = 0.2197949766366598
$ws.Range("K12").Value = 0.6231321750182417
$ws.Range("L12").Value = 0.3577104456583697
$ws.Range("N12").Value = 2.110971593708496
$ws.Range("O12").Value = 4.818874192504381
$ws.Range("B13").Value = 1.104455759757059
$ws.Range("C13").Value = 0.1383057818947293
$ws.Range("D13").Value = 0.1641119276587517
$ws.Range("F13").Value = 1.888428614489825
$ws.Range("G13").Value = 0.002494661891636437
$ws.Range("J13").Value = 0.2198306477976928
$ws.Range("K13").Value = 0.6207779293614237
$ws.Range("L13").Value = 0.3573257470126379
$ws.Range("N13").Value = 2.112058013199761
$ws.Range("O13").Value = 4.819511838332744
$ws.Range("B14").Value = 1.095533656865229
$ws.Range("C14").Value = 0.1381844674884078
$ws.Range("D14").Value = 0.1635778733361946
$ws.Range("F14").Value = 1.888620213882973
$ws.Range("G14").Value = 0.002495091602835175
$ws.Range("J14").Value = 0.2199488521550634
$ws.Range("K14").Value = 0.6131029341585474
$ws.Range("L14").Value = 0.3560744473702471
$ws.Range("N14").Value = 2.115617122277958
$ws.Range("O14").Value = 4.821644524419639
$ws.Range("B15").Value = 1.090075045031966
$ws.Range("C15").Value = 0.1381102077943623
$ws.Range("D15").Value = 0.1632518065832187
$ws.Range("F15").Value = 1.888753268991351
$ws.Range("G15").Value = 0.002495356387490213
$ws.Range("J15").Value = 0.2200226954283693
$ws.Range("K15").Value = 0.6084038723286938
$ws.Range("L15").Value = 0.3553105180923666
$ws.Range("N15").Value = 2.117809519141558
$ws.Range("O15").Value = 4.822991730101251
$ws.Range("B16").Value = 1.05889315885355
$ws.Range("C16").Value = 0.1376854174838087
$ws.Range("D16").Value = 0.1613995076339592
$ws.Range("F16").Value = 1.889756468293655
$ws.Range("G16").Value = 0.002496897928486673
$ws.Range("J16").Value = 0.220467947934015
$ws.Range("K16").Value = 0.5815079919482855
$ws.Range("L16").Value = 0.3509716038168307
$ws.Range("N16").Value = 2.130562385360406
$ws.Range("O16").Value = 4.831338752821353
$ws.Range("B17").Value = 1.039853260554622
$ws.Range("C17").Value = 0.1374255161716462
$ws.Range("D17").Value = 0.1602777521586489
$ws.Range("F17").Value = 1.890586859073508
$ws.Range("G17").Value = 0.002497865223709285
$ws.Range("J17").Value = 0.2207608259240779
$ws.Range("K17").Value = 0.5650377605573453
$ws.Range("L17").Value = 0.3483446115281623
$ws.Range("N17").Value = 2.138554607175196
$ws.Range("O17").Value = 4.837018892867576
$ws.Range("B18").Value = 1.028935045483962
$ws.Range("C18").Value = 0.137276288283033
$ws.Range("D18").Value = 0.1596379273807429
$ws.Range("F18").Value = 1.891143541234598
$ws.Range("G18").Value = 0.002498429545450476
$ws.Range("J18").Value = 0.2209365409423825
$ws.Range("K18").Value = 0.5555755582950894
$ws.Range("L18").Value = 0.3468464680247934
$ws.Range("N18").Value = 2.143213534008598
$ws.Range("O18").Value = 4.840491684164931
$ws.Range("B19").Value = 1.025244022999345
$ws.Range("C19").Value = 0.1372258075432597
$ws.Range("D19").Value = 0.1594222199655349
$ws.Range("F19").Value = 1.891345609725917
$ws.Range("G19").Value = 0.00249862198382766
$ws.Range("J19").Value = 0.2209972826480282
$ws.Range("K19").Value = 0.552373733270116
$ws.Range("L19").Value = 0.3463414309436104
$ws.Range("N19").Value = 2.144801616960087
$ws.Range("O19").Value = 4.841702855180472
$ws.Range("B20").Value = 1.041876677477802
$ws.Range("C20").Value = 0.1374531562778856
$ws.Range("D20").Value = 0.1603966087686644
$ws.Range("F20").Value = 1.890490281286787
$ws.Range("G20").Value = 0.002497761430079198
$ws.Range("J20").Value = 0.220728897441699
$ws.Range("K20").Value = 0.5667899065877577
$ws.Range("L20").Value = 0.3486229322452914
$ws.Range("N20").Value = 2.137697403413187
$ws.Range("O20").Value = 4.836392942891337
$ws.Range("B21").Value = 1.098152955393033
$ws.Range("C21").Value = 0.1382200903139648
$ws.Range("D21").Value = 0.1637345179094467
$ws.Range("F21").Value = 1.888560666050765
$ws.Range("G21").Value = 0.002494965058140368
$ws.Range("J21").Value = 0.2199138326871157
$ws.Range("K21").Value = 0.6153568308951378
$ws.Range("L21").Value = 0.3564414578620756
$ws.Range("N21").Value = 2.114569154044109
$ws.Range("O21").Value = 4.821009585328881
$ws.Range("B22").Value = 1.135221412166175
$ws.Range("C22").Value = 0.1387235287336779
$ws.Range("D22").Value = 0.1659634075947878
$ws.Range("F22").Value = 1.888002197392836
$ws.Range("G22").Value = 0.002493207968108554
$ws.Range("J22").Value = 0.2194456011615458
$ws.Range("K22").Value = 0.647192374386691
$ws.Range("L22").Value = 0.3616645858697183
$ws.Range("N22").Value = 2.100005796080551
$ws.Range("O22").Value = 4.812785375203191
$ws.Range("B23").Value = 1.115411132467329
$ws.Range("C23").Value = 0.1384546382623455
$ws.Range("D23").Value = 0.1647694981845262
$ws.Range("F23").Value = 1.888236038320258
$ws.Range("G23").Value = 0.00249413932738141
$ws.Range("J23").Value = 0.2196896158544632
$ws.Range("K23").Value = 0.630192717946386
$ws.Range("L23").Value = 0.3588665887298106
$ws.Range("N23").Value = 2.107727995372398
$ws.Range("O23").Value = 4.817007481800715
$ws.Range("B24").Value = 1.040961803171712
$ws.Range("C24").Value = 0.1374406595852236
$ws.Range("D24").Value = 0.1603428578438297
$ws.Range("F24").Value = 1.89053369713821
$ws.Range("G24").Value = 0.002497808329500603
$ws.Range("J24").Value = 0.2207433094575357
$ws.Range("K24").Value = 0.5659977402875143
$ws.Range("L24").Value = 0.34849706559514
$ws.Range("N24").Value = 2.138084745588906
$ws.Range("O24").Value = 4.836675289422573
$ws.Range("B25").Value = 0.9619219090360502
$ws.Range("C25").Value = 0.1363569435752723
$ws.Range("D25").Value = 0.1557739632100521
$ws.Range("F25").Value = 1.896035631115055
$ws.Range("G25").Value = 0.002502070476528096
$ws.Range("J25").Value = 0.2221575365908031
$ws.Range("K25").Value = 0.4971772499131077
$ws.Range("L25").Value = 0.3378030318879723
$ws.Range("N25").Value = 2.173200326689597
$ws.Range("O25").Value = 4.865756760904162
